$wb = $excel.ActiveWorkbook

# 1. Rename the "Include from LOINC" sheet to "Include #0"
$includeSheet = $wb.Worksheets.Item("Include from LOINC")
$includeSheet.Name = "Include #0"

# 2. Update the Metadata sheet
$ws = $wb.Worksheets.Item("Metadata")

# Insert a new row before row 11 (currently "Description"), shifting the
# existing rows 11-14 down to 12-15, and populate it with the new
# "Jurisdiction" property (empty value), matching the formatting of the
# other data rows.
$ws.Range("A11:B11").Insert(-4121)
$ws.Range("A10:B10").Copy()
$ws.Range("A11:B11").PasteSpecial(-4122)
$ws.Cells.Item(11, 1).Value = "Jurisdiction"
$ws.Cells.Item(11, 2).Value = ""

# Update the Date value (row 8, column B)
$ws.Cells.Item(8, 2).Value = "2024-09-17T19:55:11+00:00"
